$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right=4 (was 5), Wrong=-2 (was -1)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right=88 (was 110), Max label text updated to "88 / 112"
$ws.Range("B12").Value = 88
$ws.Range("E12").Value = "88 / 112"
